# Update odds values on Sheet1 (rows 2, 3, 7, 9, 11) to match the
# refreshed FlashScore odds snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 1.95
$ws.Range("I2").Value = 4.5
$ws.Range("Y2").Value = 10
$ws.Range("Z2").Value = 17
$ws.Range("BA2").Value = 151

# Row 3
$ws.Range("G3").Value = 1.48
$ws.Range("H3").Value = 4.33
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 2.05
$ws.Range("Q3").Value = 2.1
$ws.Range("R3").Value = 1.7
$ws.Range("W3").Value = 5
$ws.Range("AB3").Value = 41
$ws.Range("AD3").Value = 8.5
$ws.Range("AG3").Value = 13
$ws.Range("AL3").Value = 67
$ws.Range("AN3").Value = 3.2
$ws.Range("AQ3").Value = 23
$ws.Range("AW3").Value = 8.5
$ws.Range("AY3").Value = 51
$ws.Range("BC3").Value = 151

# Row 7
$ws.Range("G7").Value = 1.36
$ws.Range("H7").Value = 4.75
$ws.Range("I7").Value = 8.5
$ws.Range("N7").Value = 13
$ws.Range("O7").Value = 1.25
$ws.Range("P7").Value = 3.75
$ws.Range("Q7").Value = 1.8
$ws.Range("R7").Value = 2
$ws.Range("W7").Value = 6.5
$ws.Range("Y7").Value = 8.5
$ws.Range("AG7").Value = 19
$ws.Range("AI7").Value = 23
$ws.Range("AJ7").Value = 101
$ws.Range("AW7").Value = 9
$ws.Range("BA7").Value = 201

# Row 9
$ws.Range("G9").Value = 2.55
$ws.Range("H9").Value = 2.88
$ws.Range("I9").Value = 3.1
$ws.Range("O9").Value = 1.53
$ws.Range("P9").Value = 2.38
$ws.Range("Q9").Value = 2.7
$ws.Range("R9").Value = 1.44

# Row 11
$ws.Range("N11").Value = 8.5
$ws.Range("AQ11").Value = 34
$ws.Range("AZ11").Value = 101
